$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 506, shifting existing rows 506-577 down to 507-578
$ws.Rows.Item(506).Insert()

# Populate the new row 506 with a new weekly price observation (same
# market/category template as the surrounding rows, new date + prices).
$ws.Range("A506").Value = 3
$ws.Range("B506").Value = "Femacal de La Calera"
$ws.Range("C506").Value = "Coquimbo"
$ws.Range("D506").Value = 45131
$ws.Range("E506").Value = 5
$ws.Range("F506").Value = 100112012
$ws.Range("G506").Value = "Espinaca"
$ws.Range("H506").Value = "Sin especificar"
$ws.Range("I506").Value = "Primera"
$ws.Range("J506").Value = 110
$ws.Range("K506").Value = 4500
$ws.Range("L506").Value = 4500
$ws.Range("M506").Value = 4500
$ws.Range("N506").Value = "$/docena de atados (3 kilos)"
$ws.Range("O506").Value = "Provincia de Quillota"
$ws.Range("P506").Value = 1500
$ws.Range("Q506").Value = 3
$ws.Range("R506").Value = "Hortaliza"
